$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = 'Other found locations'
$ws.Range("E2").Value = '[Simon D%Harding%NULL%1,   Joanna L%Sharman%NULL%2,   Joanna L%Sharman%NULL%0,   Elena%Faccenda%NULL%1,   Chris%Southan%NULL%1,   Adam J%Pawson%NULL%1,   Sam%Ireland%NULL%1,   Alasdair J G%Gray%NULL%1,   Liam%Bruce%NULL%1,   Stephen P H%Alexander%NULL%1,   Stephen%Anderton%NULL%1,   Clare%Bryant%NULL%1,   Anthony P%Davenport%NULL%1,   Christian%Doerig%NULL%1,   Doriano%Fabbro%NULL%1,   Francesca%Levi-Schaffer%NULL%1,   Michael%Spedding%NULL%1,   Jamie A%Davies%jamie.davies@ed.ac.uk%1,   NULL%NULL%NULL%0]'
$ws.Range("I2").Value = '_PMC'
$ws.Range("E3").Value = '[ Ourania S.%Kotsiou%null%2,    Sotirios G.%Zarogiannis%null%1,    Konstantinos I.%Gourgoulianis%null%1,  Ourania S.%Kotsiou%null%0,  Sotirios G.%Zarogiannis%null%1,  Konstantinos I.%Gourgoulianis%null%1]'
$ws.Range("I3").Value = ''
$ws.Range("E4").Value = '[ Jonathan%Messika%null%2,    Benjamin%Sztrymf%null%1,    Fabrice%Bertrand%null%1,    Typhaine%Billard-Pomares%null%1,    Guilene%Barnaud%null%1,    Catherine%Branger%null%1,    Didier%Dreyfuss%null%1,    Jean-Damien%Ricard%null%1,  Jonathan%Messika%null%0,  Benjamin%Sztrymf%null%1,  Fabrice%Bertrand%null%1,  Typhaine%Billard-Pomares%null%1,  Guilene%Barnaud%null%1,  Catherine%Branger%null%1,  Didier%Dreyfuss%null%1,  Jean-Damien%Ricard%null%1]'
$ws.Range("I4").Value = ''
$ws.Range("E5").Value = '[Carl%Llor%carles.llor@urv.cat%1,   Ana%Moragas%amoragasm@meditex.es%1,   Carolina%Bayona%cbayona.tarte.ics@gencat.cat%1,   Rosa%Morros%rmorros@idiapjgol.org%1,   Helena%Pera%hpera@idiapjgol.org%1,   Josep M%Cots%23465jcy@comb.cat%1,   Yvonne%Fernández%yfernandez.tarte.ics@gencat.cat%1,   Marc%Miravitlles%marcm@separ.es%1,   Albert%Boada%aboadav.bcn.ics@gencat.cat%1]'
$ws.Range("I5").Value = '_PMC_Springer'
$ws.Range("I6").Value = ''
$ws.Range("E7").Value = '[ G.%Voiriot%null%2,    M.%Chalumeau%null%1,    J.%Messika%null%1,    D.%Basille%null%1,    B.%Philippe%null%1,    J.-D.%Ricard%null%1,    C.%Andrejak%null%1,    V.%Jounieaux%null%1,    O.%Sanchez%null%1,    M.%Fartoukh%null%1,  G.%Voiriot%null%0,  M.%Chalumeau%null%1,  J.%Messika%null%1,  D.%Basille%null%1,  B.%Philippe%null%1,  J.-D.%Ricard%null%1,  C.%Andrejak%null%1,  V.%Jounieaux%null%1,  O.%Sanchez%null%1,  M.%Fartoukh%null%1]'
$ws.Range("I7").Value = ''
$ws.Range("E8").Value = '[ Damien%Basille%null%2,    Nathalie%Plouvier%null%1,    Charlotte%Trouve%null%1,    Pierre%Duhaut%null%1,    Claire%Andrejak%null%1,    Vincent%Jounieaux%null%1,  Damien%Basille%null%0,  Nathalie%Plouvier%null%1,  Charlotte%Trouve%null%1,  Pierre%Duhaut%null%1,  Claire%Andrejak%null%1,  Vincent%Jounieaux%null%1]'
$ws.Range("I8").Value = ''
$ws.Range("I9").Value = ''
